$d = $word.ActiveDocument

# wdColorAutomatic / explicit black
$wdColorAutomatic = -16777216
$wdColorBlack = 0

# ---------------------------------------------------------------
# Change 1: "...well-forged as Vienna swords, in Dauphiné. The..."
#        -> "...well-forged as swords from Vienne in Dauphiné. The..."
# ---------------------------------------------------------------

$d.Content.Find.Execute("Vienna swords, in Dauphiné", $true, $false, $false, $false, $false, $true, 1, $false, "swords from Vienne in Dauphiné", 2) | Out-Null

$rng1 = $d.Content
$rng1.Find.Execute("swords from Vienne in Dauphiné", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$base1 = $rng1.Start

# "swords from " (no explicit color)
$d.Range($base1, $base1 + 12).Font.Color = $wdColorAutomatic
# "Vienn" (explicit black)
$d.Range($base1 + 12, $base1 + 17).Font.Color = $wdColorBlack
# "e in " (no explicit color)
$d.Range($base1 + 17, $base1 + 22).Font.Color = $wdColorAutomatic
# "Dauphiné..." onward keeps its existing explicit black color (untouched)

# ---------------------------------------------------------------
# Change 2: "<ab>Vienna blades cost commonly Xviii or 19 lb.</ab>"
#        -> "<ab>Blades from Vienne commonly cost xviii or 19 lb.</ab>"
# ---------------------------------------------------------------

$d.Content.Find.Execute("Vienna blades cost commonly Xviii or 19 lb.", $true, $false, $false, $false, $false, $true, 1, $false, "Blades from Vienne commonly cost xviii or 19 lb.", 2) | Out-Null

$rng2 = $d.Content
$rng2.Find.Execute("<ab>Blades from Vienne commonly cost xviii or 19 lb.</ab>", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$base2 = $rng2.Start

# "<ab>" (explicit black) -- left untouched, already black from original run

# "B" (no explicit color)
$d.Range($base2 + 4, $base2 + 5).Font.Color = $wdColorAutomatic
# "lades from Vienne" (explicit black)
$d.Range($base2 + 5, $base2 + 22).Font.Color = $wdColorBlack
# " commonly " (no explicit color)
$d.Range($base2 + 22, $base2 + 32).Font.Color = $wdColorAutomatic
# "cost " (explicit black)
$d.Range($base2 + 32, $base2 + 37).Font.Color = $wdColorBlack
# "x" (no explicit color)
$d.Range($base2 + 37, $base2 + 38).Font.Color = $wdColorAutomatic
# "viii or 19 lb.</ab>" (explicit black)
$d.Range($base2 + 38, $base2 + 57).Font.Color = $wdColorBlack
